# Auto-generated edits applying the diff to Marilith_Profits workbook
$wb = $excel.ActiveWorkbook

$edits = @(
    @("ALC", "H2", 0)
    @("ALC", "I2", 0)
    @("ALC", "J2", 0)
    @("ALC", "K2", 0)
    @("ALC", "L2", 0)
    @("ALC", "M2", $null)
    @("ALC", "N2", $null)
    @("ALC", "H13", 39999)
    @("ALC", "I13", 39999)
    @("ALC", "K13", 39999)
    @("ALC", "M13", -39830)
    @("ALC", "H20", 0)
    @("ALC", "I20", 0)
    @("ALC", "K20", 0)
    @("ALC", "M20", $null)
    @("ALC", "H33", 263.47058)
    @("ALC", "I33", 192.78572)
    @("ALC", "J33", 593.3333)
    @("ALC", "K33", 192.78572)
    @("ALC", "L33", 593.3333)
    @("ALC", "M33", 36.21428)
    @("ALC", "N33", -1051.3333)
    @("ALC", "H35", 0)
    @("ALC", "I35", 0)
    @("ALC", "K35", 0)
    @("ALC", "M35", $null)
    @("ALC", "H38", 346.22223)
    @("ALC", "J38", 1000)
    @("ALC", "L38", 3000)
    @("ALC", "N38", -3744)
    @("ALC", "H51", 3000)
    @("ALC", "I51", 2500)
    @("ALC", "J51", 3166.6667)
    @("ALC", "K51", 2500)
    @("ALC", "L51", 3166.6667)
    @("ALC", "M51", -2016)
    @("ALC", "N51", -4134.6667)
    @("ALC", "H53", 1120.2)
    @("ALC", "I53", 1233)
    @("ALC", "J53", 951)
    @("ALC", "K53", 1233)
    @("ALC", "L53", 951)
    @("ALC", "M53", -596)
    @("ALC", "N53", -2225)
    @("ALC", "H62", 7185.625)
    @("ALC", "I62", 6497.2)
    @("ALC", "K62", 6497.2)
    @("ALC", "M62", -5873.2)
    @("ALC", "H65", 7185.625)
    @("ALC", "I65", 6497.2)
    @("ALC", "K65", 32486)
    @("ALC", "M65", -29366)
    @("ALC", "H69", 2332.6667)
    @("ALC", "I69", 2500)
    @("ALC", "J69", 1998)
    @("ALC", "K69", 7500)
    @("ALC", "L69", 5994)
    @("ALC", "M69", -6626)
    @("ALC", "N69", -7742)
    @("ALC", "H72", 2332.6667)
    @("ALC", "I72", 2500)
    @("ALC", "J72", 1998)
    @("ALC", "K72", 22500)
    @("ALC", "L72", 17982)
    @("ALC", "M72", -18132)
    @("ALC", "N72", -26718)
    @("ALC", "H93", 0)
    @("ALC", "J93", 0)
    @("ALC", "L93", 0)
    @("ALC", "N93", $null)
    @("ALC", "H98", 1152.4375)
    @("ALC", "I98", 1166.6666)
    @("ALC", "J98", 1109.75)
    @("ALC", "K98", 1166.6666)
    @("ALC", "L98", 1109.75)
    @("ALC", "M98", 331.3334)
    @("ALC", "N98", -4105.75)
    @("ALC", "H113", 6763.7144)
    @("ALC", "I113", 5188)
    @("ALC", "J113", 9600)
    @("ALC", "K113", 5188)
    @("ALC", "L113", 9600)
    @("ALC", "M113", -1934)
    @("ALC", "N113", -16108)
    @("ALC", "H115", 93.666664)
    @("ALC", "I115", 93.666664)
    @("ALC", "K115", 280.999992)
    @("ALC", "M115", 1286.000008)
    @("ALC", "H122", 1152.4375)
    @("ALC", "I122", 1166.6666)
    @("ALC", "J122", 1109.75)
    @("ALC", "K122", 3499.9998)
    @("ALC", "L122", 3329.25)
    @("ALC", "M122", -1049.9998)
    @("ALC", "N122", -8229.25)
    @("ALC", "H135", 933)
    @("ALC", "I135", 933)
    @("ALC", "K135", 8397)
    @("ALC", "M135", -5862)
    @("ARM", "H32", 1506.6316)
    @("ARM", "I32", 1292.9143)
    @("ARM", "J32", 4000)
    @("ARM", "K32", 1292.9143)
    @("ARM", "L32", 4000)
    @("ARM", "M32", -1005.9143)
    @("ARM", "N32", -4574)
    @("ARM", "H102", 2925)
    @("ARM", "J102", 3000)
    @("ARM", "L102", 3000)
    @("ARM", "N102", -6244)
    @("ARM", "H132", 405.41666)
    @("ARM", "I132", 405.41666)
    @("ARM", "K132", 1216.24998)
    @("ARM", "M132", 1313.75002)
    @("BSM", "H20", 1070.0834)
    @("BSM", "I20", 879.1429000000001)
    @("BSM", "K20", 879.1429000000001)
    @("BSM", "M20", -632.1429000000001)
    @("BSM", "H22", 147.5)
    @("BSM", "I22", 145)
    @("BSM", "J22", 150)
    @("BSM", "K22", 145)
    @("BSM", "L22", 150)
    @("BSM", "M22", 28)
    @("BSM", "N22", -496)
    @("BSM", "H107", 1297.3334)
    @("BSM", "I107", 1297.3334)
    @("BSM", "J107", 0)
    @("BSM", "K107", 1297.3334)
    @("BSM", "L107", 0)
    @("BSM", "M107", 622.6666)
    @("BSM", "N107", $null)
    @("BSM", "H134", 7015.579)
    @("BSM", "I134", 7076.778)
    @("BSM", "K134", 21230.334)
    @("BSM", "M134", -18695.334)
    @("CRP", "H107", 804.0714)
    @("CRP", "I107", 718.9091)
    @("CRP", "J107", 1116.3334)
    @("CRP", "K107", 718.9091)
    @("CRP", "L107", 1116.3334)
    @("CRP", "M107", 1201.0909)
    @("CRP", "N107", -4956.3334)
    @("CUL", "H54", 995)
    @("CUL", "J54", 995)
    @("CUL", "L54", 2985)
    @("CUL", "N54", -4103)
    @("CUL", "H134", 340666.34)
    @("CUL", "J134", 11000)
    @("CUL", "L134", 33000)
    @("CUL", "N134", -43140)
    @("GSM", "H2", 161.64285)
    @("GSM", "I2", 202.77777)
    @("GSM", "J2", 87.59999999999999)
    @("GSM", "K2", 202.77777)
    @("GSM", "L2", 87.59999999999999)
    @("GSM", "M2", -89.77777)
    @("GSM", "N2", -313.6)
    @("GSM", "H19", 11856.714)
    @("GSM", "I19", 8999.666999999999)
    @("GSM", "K19", 8999.666999999999)
    @("GSM", "M19", -8711.666999999999)
    @("GSM", "H107", 115.55556)
    @("GSM", "I107", 121.42857)
    @("GSM", "K107", 121.42857)
    @("GSM", "M107", 1798.57143)
    @("GSM", "H113", 2000)
    @("GSM", "I113", 2000)
    @("GSM", "K113", 2000)
    @("GSM", "M113", 170)
    @("LTW", "H22", 690)
    @("LTW", "I22", 500)
    @("LTW", "J22", 785)
    @("LTW", "K22", 500)
    @("LTW", "L22", 785)
    @("LTW", "M22", -205)
    @("LTW", "N22", -1375)
    @("LTW", "H27", 690)
    @("LTW", "I27", 500)
    @("LTW", "J27", 785)
    @("LTW", "K27", 500)
    @("LTW", "L27", 785)
    @("LTW", "M27", -393)
    @("LTW", "N27", -999)
    @("LTW", "H40", 7997.4)
    @("LTW", "I40", 7496.75)
    @("LTW", "K40", 7496.75)
    @("LTW", "M40", -7360.75)
    @("LTW", "H55", 235.45454)
    @("LTW", "I55", 215.83333)
    @("LTW", "J55", 259)
    @("LTW", "K55", 215.83333)
    @("LTW", "L55", 259)
    @("LTW", "M55", -42.83332999999999)
    @("LTW", "N55", -605)
    @("LTW", "H122", 3083.7273)
    @("LTW", "I122", 3083.7273)
    @("LTW", "K122", 9251.1819)
    @("LTW", "M122", -6801.1819)
    @("LTW", "H132", 0)
    @("LTW", "I132", 0)
    @("LTW", "K132", 0)
    @("LTW", "M132", $null)
    @("LTW", "H6", 6000)
    @("LTW", "J6", 6000)
    @("LTW", "L6", 6000)
    @("LTW", "N6", -6230)
    @("WVR", "H126", 26783.357)
    @("WVR", "I126", 19496.8)
    @("WVR", "J126", 44999.75)
    @("WVR", "K126", 58490.39999999999)
    @("WVR", "L126", 134999.25)
    @("WVR", "M126", -56020.39999999999)
    @("WVR", "N126", -139939.25)
    @("WVR", "H132", 2341.6667)
    @("WVR", "I132", 2183.3333)
    @("WVR", "J132", 2500)
    @("WVR", "K132", 6549.999899999999)
    @("WVR", "L132", 7500)
    @("WVR", "M132", -4019.999899999999)
    @("WVR", "N132", -12560)
)

foreach ($edit in $edits) {
    $ws = $wb.Worksheets.Item($edit[0])
    $rng = $ws.Range($edit[1])
    if ($edit[2] -eq $null) {
        $rng.ClearContents()
    } else {
        $rng.Value = $edit[2]
    }
}
